$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = "'43.152.60"
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = '  -1.51%  '
$ws.Range('D3').Value = "'2.266.67"
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = '  -2.11%  '
$ws.Range('E4').Value = '  -0.18%  '
$ws.Range('D5').Value = "'110.75"
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -1.52%  '
$ws.Range('D6').Value = "'263.75"
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -2.54%  '
$ws.Range('D7').Value = "'0.617"
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  -0.58%  '
$ws.Range('E8').Value = '  +0.10%  '
$ws.Range('D9').Value = "'0.601"
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  -3.67%  '
$ws.Range('D10').Value = "'47.23"
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  -0.12%  '
$ws.Range('D11').Value = "'0.0928"
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  -1.63%  '
$ws.Range('D12').Value = "'8.74"
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  -1.84%  '
$ws.Range('D13').Value = "'0.108"
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  +0.32%  '
$ws.Range('D14').Value = "'15.37"
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  -2.57%  '
$ws.Range('D15').Value = "'2.605.83"
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  -2.22%  '
$ws.Range('D16').Value = "'0.848"
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  -1.67%  '
$ws.Range('D17').Value = "'2.266.02"
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  -2.35%  '
$ws.Range('D18').Value = "'43.006.73"
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  -2.01%  '
$ws.Range('E19').Value = '  -2.43%  '
$ws.Range('D20').Value = "'6.79"
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  +1.41%  '
$ws.Range('D21').Value = "'70.97"
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -2.22%  '
$ws.Range('D22').Value = "'2.50"
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -0.25%  '
$ws.Range('D23').Value = "'230.72"
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  -1.51%  '
$ws.Range('D24').Value = "'9.58"
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +0.62%  '
$ws.Range('D25').Value = "'2.85"
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  -1.35%  '
$ws.Range('D26').Value = "'0.999"
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  +0.03%  '
$ws.Range('D27').Value = "'11.25"
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -2.67%  '
$ws.Range('E28').Value = '  -1.30%  '
$ws.Range('D29').Value = "'40.18"
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  -5.60%  '
$ws.Range('E30').Value = '  -1.91%  '
$ws.Range('D31').Value = "'3.27"
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  -4.87%  '
$ws.Range('D32').Value = "'171.19"
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  -3.80%  '
$ws.Range('D33').Value = "'21.24"
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  -3.16%  '
$ws.Range('D34').Value = "'0.0898"
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  -3.95%  '
$ws.Range('E35').Value = '  +1.15%  '
$ws.Range('E36').Value = '  -0.67%  '
$ws.Range('D37').Value = "'4.64"
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  -3.11%  '
$ws.Range('E38').Value = '  -2.94%  '
$ws.Range('E39').Value = '  -6.71%  '
$ws.Range('E40').Value = '  -4.30%  '
$ws.Range('E41').Value = '  +8.30%  '
$ws.Range('D42').Value = "'75.78"
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  +8.88%  '
$ws.Range('D43').Value = "'13.79"
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  +6.70%  '
$ws.Range('D44').Value = "'0.234"
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -5.04%  '
$ws.Range('D45').Value = "'6.04"
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  +4.77%  '
$ws.Range('E46').Value = '  -0.15%  '
$ws.Range('E47').Value = '  -3.02%  '
$ws.Range('E48').Value = '  -2.60%  '
$ws.Range('D49').Value = "'0.0991"
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  -1.83%  '
$ws.Range('D50').Value = "'1.24"
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +1.18%  '
$ws.Range('D51').Value = "'100.40"
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +0.30%  '
